$d = $word.ActiveDocument
$bullet = [char]0x2022

# ============================================================
# Change 1: Condense the three CORE COMPETENCIES detail
# paragraphs (Survey Methodology..., Redistricting...,
# Data Analysis...) into a single summary paragraph.
# ============================================================
$surveyPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Survey Methodology & Research Design:")) {
        $surveyPara = $p
        break
    }
}

if ($surveyPara -ne $null) {
    $condensed = "Survey Methodology & Research Design $bullet Redistricting & Geospatial Analysis $bullet Data Analysis & Visualization"
    $surveyPara.Range.Text = $condensed

    # The paragraph(s) that used to follow (Redistricting..., Data Analysis...)
    # are now immediately after $surveyPara; delete their paragraph ranges
    # (not including the trailing paragraph mark) so they collapse away
    # without leaving blank paragraphs behind.
    $next1 = $surveyPara.Next()
    if ($next1 -ne $null -and $next1.Range.Text.StartsWith("Redistricting & Geospatial Analysis:")) {
        $next1.Range.Delete()
    }
    $next2 = $surveyPara.Next()
    if ($next2 -ne $null -and $next2.Range.Text.StartsWith("Data Analysis & Visualization:")) {
        $next2.Range.Delete()
    }
}

# ============================================================
# Change 2: Insert a new "TECHNICAL SKILLS" section (heading +
# three detail paragraphs) right before the closing
# "For a more detailed..." paragraph.
# ============================================================
$anchorPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("For a more detailed, full description")) {
        $anchorPara = $p
        break
    }
}

if ($anchorPara -ne $null) {
    $prevPara = $anchorPara.Previous()

    $newBlock = "`r" + "TECHNICAL SKILLS" `
        + "`r" + "SURVEY METHODOLOGY & RESEARCH DESIGN Survey Design and Questionnaire Development for Political and Market Research; Sampling Methodology and Statistical Analysis (R, SPSS, Stata, OSCAR); Random Device Engagement (RDE), Text Message, Web Panel, and Live Telephone Calling; Expert Testimony and Consultation on Research Methodology" `
        + "`r" + "REDISTRICTING & GEOSPATIAL ANALYSIS Redistricting Software Development and Boundary Estimation Systems; Geospatial Analysis; Choropleths and Hexagonal Grid Maps for Demographic Visualization; Court Case Analysis and Expert Testimony for Redistricting" `
        + "`r" + "DATA ANALYSIS & VISUALIZATION Advanced Statistical Modeling and Analysis (Regression, Clustering, Segmentation); Data Visualization; Consumer Behavior Analysis and Market Segmentation; Multi-million Dollar Research Project Management"

    $prevPara.Range.InsertAfter($newBlock)

    $headingPara = $prevPara.Next()
    $headingPara.Style = "Heading 2"
}

Write-Host "Done. Final paragraph count: " $d.Paragraphs.Count
